$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Sheet1" to "Active Campaign"
$ws.Name = "Active Campaign"

# Fix the typo in the campaign vanity URL (evergeen -> evergreen)
$ws.Range("A2").Value = "https://www.crepeerase.com/campaign/classic-evergreen-refresh"

# Remove the trailing "All orders are subject to applicable sales tax." paragraph
# from the Post Purchase Cart Language text in C43
$ws.Range("C43").Value = "Three months after your first order is shipped, and then every three months thereafter, you will be sent a new full size supply of Crepe Erase. Each shipment will be charged to the card you provide today, in three monthly payments at the low price of `$39.95 plus `$2.99 for shipping and handling per month, unless you call to cancel.There is no commitment and no minimum to buy."

# Row 43 is shorter now that the text is shorter - reduce its custom height
$ws.Rows.Item(43).RowHeight = 122.25

# Update the sheet view: scroll so row 28 is at the top, and select cell B43
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
[void]$ws.Range("B43").Select()

Write-Host "Edits applied"
